$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The category previously labelled "Baby Products" (row 19, column B) is
# renamed to "Baby Images" now that the UI pulls it from Firebase.
$ws.Range("B19").Value = "Baby Images"

# Reflect the saved window/selection state: the user had scrolled back up
# and landed on C10 before saving.
$ws.Range("C10").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
